$d = $word.ActiveDocument

# Locate the insertion point: right after the 3rd occurrence of "(Billede)"
# (the end of the "Tilføj_vare_udkast.jpg" bilag line), which is immediately
# before the "Bilag<tab>Kodedokumentation...(Doxygen HTML)" run that we
# need to replace/expand.
$startRng = $d.Content.Duplicate
$fStart = $startRng.Find
$fStart.Text = "(Billede)"
$fStart.Forward = $true
$fStart.Wrap = 0
$fStart.Execute() | Out-Null
$fStart.Execute() | Out-Null
$fStart.Execute() | Out-Null
$insertStart = $startRng.End

# Locate the end of the region to replace: the end of the (only) "(Doxygen HTML)"
# occurrence, plus the (invisible) _GoBack bookmark and the line break that
# immediately follow it (these get reconstructed inside the replacement XML).
$endRng = $d.Content.Duplicate
$fEnd = $endRng.Find
$fEnd.Text = "(Doxygen HTML)"
$fEnd.Forward = $true
$fEnd.Wrap = 0
$fEnd.Execute() | Out-Null
$insertEnd = $endRng.End + 1

$target = $d.Range($insertStart, $insertEnd)

$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:br/><w:t>Bilag</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>Kode_FridgeApp</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>(VS2013 solution)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:br/><w:t>Bilag</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>Kode_</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Web</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:i/></w:rPr><w:t>App</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>(VS2013 solution)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:br/><w:t>Bilag</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>Kodedokumentation_FrigeApp</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>(Doxygen HTML)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:br/><w:t>Bilag</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>Kodedokumentation_</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Web</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>App</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t>(Doxygen HTML)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newXml)
